# Fix 0 hospitalisation error for India
# Adds two new lookup columns (is_india / type) to the wfa dictionary sheet,
# flagging every existing row as an India row of "character" type.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New column headers
$ws.Range("D1").Value = "is_india"
$ws.Range("E1").Value = "type"

# Populate every existing data row (2-24) with the new flag/type values
for ($r = 2; $r -le 24; $r++) {
    $ws.Cells.Item($r, 4).Value = 1
    $ws.Cells.Item($r, 5).Value = "character"
}

# Match the centered alignment already used by the other helper columns
$ws.Range("D2:D24").HorizontalAlignment = -4108
$ws.Range("D2:D24").VerticalAlignment = -4108

# Move the active selection to where the author left off editing
$ws.Range("H15").Select() | Out-Null
